# Auto-generated data array: header row + 65 component rows (Designator, Mid X, Mid Y, Layer, Rotation)
$data = @(
  @("Designator","Mid X","Mid Y","Layer","Rotation"),
  @("C1",-16.446599,-38.576599,"Top",90.0),
  @("C10",0.32,-32.63,"Top",90.0),
  @("C11",3.81,-32.53,"Top",90.0),
  @("C12",23.97,-22.23,"Top",90.0),
  @("C13",17.79,-23.31,"Top",180.0),
  @("C14",13.68,-23.37,"Top",180.0),
  @("C15",14.94,-7.4,"Top",180.0),
  @("C16",18.3734,3.5234,"Top",180.0),
  @("C17",-18.823399,3.5234,"Top",180.0),
  @("C18",-18.823399,-1.526599,"Top",180.0),
  @("C19",-10.693999,4.797,"Top",0.0),
  @("C1A",-0.65,43.83,"Top",270.0),
  @("C2",-14.69,-31.22,"Top",0.0),
  @("C20",1.888,11.442,"Top",90.0),
  @("C21",2.688,5.142,"Top",180.0),
  @("C22",6.938,10.992,"Top",270.0),
  @("C23",23.83,-3.53,"Top",180.0),
  @("C24",23.85,-0.71,"Top",180.0),
  @("C25",24.08,-10.66,"Top",180.0),
  @("C26",24.2134,-13.806599,"Top",180.0),
  @("C2A",-10.21,51.68,"Top",90.0),
  @("C3",-9.65,24.61,"Top",0.0),
  @("C3A",10.22,51.49,"Top",90.0),
  @("C4",-16.47,-20.74,"Top",90.0),
  @("C5",8.71,24.41,"Top",0.0),
  @("C6",-19.85,-16.17,"Top",90.0),
  @("C7",-15.36,-16.29,"Top",270.0),
  @("C8",-2.22,-28.11,"Top",90.0),
  @("C9",-6.77,-32.5,"Top",90.0),
  @("D1A",-4.65,43.66,"Top",270.0),
  @("Q1",-18.15,-10.04,"Top",270.0),
  @("Q2",-11.2,-10.24,"Top",90.0),
  @("Q3",9.08,-28.21,"Top",180.0),
  @("Q4",10.796,4.142,"Top",0.0),
  @("R1",-11.82,-38.71,"Top",270.0),
  @("R10",-8.407999,19.446,"Top",180.0),
  @("R11",-2.24,-32.73,"Top",270.0),
  @("R12",8.064,-10.735999,"Top",90.0),
  @("R13",-11.956,-26.538999,"Top",180.0),
  @("R14",21.87,-22.63,"Top",270.0),
  @("R15",-11.189999,-0.231999,"Top",0.0),
  @("R16",9.31,-22.24,"Top",90.0),
  @("R17",-7.65,-25.46,"Top",90.0),
  @("R18",18.98,-7.3,"Top",180.0),
  @("R19",10.44,27.49,"Top",90.0),
  @("R1A",-2.45,43.85,"Top",90.0),
  @("R2",-18.91,-38.57,"Top",270.0),
  @("R20",3.933,18.671,"Top",270.0),
  @("R21",11.967,10.946,"Top",270.0),
  @("R22",8.941,18.671,"Top",90.0),
  @("R23",18.855,-9.960999,"Top",0.0),
  @("R24",24.09,-17.22,"Top",180.0),
  @("R25",13.51,37.87,"Top",90.0),
  @("R2A",1.08,43.86,"Top",270.0),
  @("R3",-21.78,-38.66,"Top",270.0),
  @("R3A",2.82,43.86,"Top",90.0),
  @("R4",-19.7,-31.27,"Top",0.0),
  @("R5",-19.73,-28.31,"Top",180.0),
  @("R6",-19.97,-21.14,"Top",90.0),
  @("R7",-23.34,-7.82,"Top",90.0),
  @("R8",-11.19,-5.24,"Top",0.0),
  @("R9",-11.8,-16.34,"Top",90.0),
  @("RLEDFX",-4.65,0.54,"Top",270.0),
  @("U1",-1.09,-8.49,"Top",180.0),
  @("U2",-6.161999,12.192,"Top",270.0)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = $data.Count
$colCount = 5

$arr = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
  for ($c = 0; $c -lt $colCount; $c++) {
    $arr[$r,$c] = $data[$r][$c]
  }
}

$targetRange = $ws.Range("A1").Resize($rowCount, $colCount)
$targetRange.Value = $arr

# Resize the existing table (Table1) to cover the new data extent
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E" + $rowCount))

$ws.Range("A1").Select()
